# Rename worksheet "Table Y" to "Table A1" (per commit: "update maps and code")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table Y")
$ws.Name = "Table A1"
